$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.729.04"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "3.469.90"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'578.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "'147.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.94%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.481"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").Value = "'7.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").Value = "'0.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("E11").Value = "  +4.24%  "
$ws.Range("D12").Value = "4.068.02"
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("D13").Value = "'29.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.29%  "
$ws.Range("D14").Value = "'0.129"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").Value = "3.467.40"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "62.821.01"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").Value = "'6.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.50%  "
$ws.Range("D19").Value = "'14.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.40%  "
$ws.Range("D20").Value = "'9.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.21%  "
$ws.Range("D21").Value = "'387.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "'0.561"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").Value = "'74.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "3.610.99"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").Value = "'0.179"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.06%  "
$ws.Range("D28").Value = "'7.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'8.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").Value = "'2.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "'1.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("D34").Value = "'23.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'7.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'5.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.94%  "
$ws.Range("D37").Value = "'1.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.59%  "
$ws.Range("D38").Value = "'31.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +22.76%  "
$ws.Range("D39").Value = "'169.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "3.511.69"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").Value = "'0.0763"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "'0.798"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'4.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'42.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").Value = "'1.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("E46").Value = "  +3.69%  "
$ws.Range("D47").Value = "2.606.51"
$ws.Range("E47").Value = "  +6.23%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.50%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'23.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("D50").Value = "'6.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
